$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bgn"
$ws.Range("C2").Value = "Tlr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 13.30571766666667
$ws.Range("H2").Value = 39.917153
$ws.Range("I2").Value = 0.007643519924167935
$ws.Range("J2").Value = 0.007643519924167933
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.493155
$ws.Range("N2").Value = 1.479465
$ws.Range("O2").Value = 0.2262533155038342
$ws.Range("P2").Value = 0.2262533155038342
$ws.Range("Q2").Value = 6.561781195905
$ws.Range("R2").Value = 59.056030763145
$ws.Range("S2").Value = 0.001729371724962611
$ws.Range("T2").Value = 0.001729371724962611

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bgn"
$ws.Range("C3").Value = "Tlr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 13.30571766666667
$ws.Range("H3").Value = 39.917153
$ws.Range("I3").Value = 0.007643519924167935
$ws.Range("J3").Value = 0.007643519924167933
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 1.648742666666666
$ws.Range("N3").Value = 4.946228
$ws.Range("O3").Value = 0.7564224123165462
$ws.Range("P3").Value = 0.7564224123165462
$ws.Range("Q3").Value = 21.93770442765377
$ws.Range("R3").Value = 197.439339848884
$ws.Range("S3").Value = 0.005781729779628693
$ws.Range("T3").Value = 0.005781729779628692

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Bgn"
$ws.Range("C4").Value = "Tlr2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 13.30571766666667
$ws.Range("H4").Value = 39.917153
$ws.Range("I4").Value = 0.007643519924167935
$ws.Range("J4").Value = 0.007643519924167933
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.037761
$ws.Range("N4").Value = 0.113283
$ws.Range("O4").Value = 0.01732427217961956
$ws.Range("P4").Value = 0.01732427217961956
$ws.Range("Q4").Value = 0.502437204811
$ws.Range("R4").Value = 4.521934843298999
$ws.Range("S4").Value = 0.0001324184195766304
$ws.Range("T4").Value = 0.0001324184195766303

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Bgn"
$ws.Range("C5").Value = "Tlr2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1644.738728666666
$ws.Range("H5").Value = 4934.216186
$ws.Range("I5").Value = 0.9448263940026712
$ws.Range("J5").Value = 0.9448263940026712
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.493155
$ws.Range("N5").Value = 1.479465
$ws.Range("O5").Value = 0.2262533155038342
$ws.Range("P5").Value = 0.2262533155038342
$ws.Range("Q5").Value = 811.1111277356099
$ws.Range("R5").Value = 7300.00014962049
$ws.Range("S5").Value = 0.2137701042186364
$ws.Range("T5").Value = 0.2137701042186364

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Bgn"
$ws.Range("C6").Value = "Tlr2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1644.738728666666
$ws.Range("H6").Value = 4934.216186
$ws.Range("I6").Value = 0.9448263940026712
$ws.Range("J6").Value = 0.9448263940026712
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 1.648742666666666
$ws.Range("N6").Value = 4.946228
$ws.Range("O6").Value = 0.7564224123165462
$ws.Range("P6").Value = 0.7564224123165462
$ws.Range("Q6").Value = 2711.750917471823
$ws.Range("R6").Value = 24405.75825724641
$ws.Range("S6").Value = 0.714687860171844
$ws.Range("T6").Value = 0.714687860171844

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Bgn"
$ws.Range("C7").Value = "Tlr2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1644.738728666666
$ws.Range("H7").Value = 4934.216186
$ws.Range("I7").Value = 0.9448263940026712
$ws.Range("J7").Value = 0.9448263940026712
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.037761
$ws.Range("N7").Value = 0.113283
$ws.Range("O7").Value = 0.01732427217961956
$ws.Range("P7").Value = 0.01732427217961956
$ws.Range("Q7").Value = 62.10697913318199
$ws.Range("R7").Value = 558.962812198638
$ws.Range("S7").Value = 0.01636842961219075
$ws.Range("T7").Value = 0.01636842961219075

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Bgn"
$ws.Range("C8").Value = "Tlr2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 82.73961633333333
$ws.Range("H8").Value = 248.218849
$ws.Range("I8").Value = 0.04753008607316088
$ws.Range("J8").Value = 0.04753008607316087
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.493155
$ws.Range("N8").Value = 1.479465
$ws.Range("O8").Value = 0.2262533155038342
$ws.Range("P8").Value = 0.2262533155038342
$ws.Range("Q8").Value = 40.803455492865
$ws.Range("R8").Value = 367.231099435785
$ws.Range("S8").Value = 0.01075383956023527
$ws.Range("T8").Value = 0.01075383956023526

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Bgn"
$ws.Range("C9").Value = "Tlr2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 82.73961633333333
$ws.Range("H9").Value = 248.218849
$ws.Range("I9").Value = 0.04753008607316088
$ws.Range("J9").Value = 0.04753008607316087
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 1.648742666666666
$ws.Range("N9").Value = 4.946228
$ws.Range("O9").Value = 0.7564224123165462
$ws.Range("P9").Value = 0.7564224123165462
$ws.Range("Q9").Value = 136.4163356723969
$ws.Range("R9").Value = 1227.747021051572
$ws.Range("S9").Value = 0.03595282236507343
$ws.Range("T9").Value = 0.03595282236507343

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Bgn"
$ws.Range("C10").Value = "Tlr2"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 82.73961633333333
$ws.Range("H10").Value = 248.218849
$ws.Range("I10").Value = 0.04753008607316088
$ws.Range("J10").Value = 0.04753008607316087
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.037761
$ws.Range("N10").Value = 0.113283
$ws.Range("O10").Value = 0.01732427217961956
$ws.Range("P10").Value = 0.01732427217961956
$ws.Range("Q10").Value = 3.124330652363
$ws.Range("R10").Value = 28.118975871267
$ws.Range("S10").Value = 0.0008234241478521842
$ws.Range("T10").Value = 0.0008234241478521841
